# "Generate Report for Handback" — refresh the localization-status report
# after a successful handback: the Overview / per-locale "Status" columns
# move from "Ready for handoff" to "Handed back: in sync with en-US", the
# handback timestamps + error details for the already-synced files are
# updated/cleared, and the columns that changed content get re-sized to
# fit their new text.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for both rows.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

# Columns widened to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

# New handback datetimes recorded for this run.
$wsZh.Range("K2").Value = "2016-11-29 02:46:02"
$wsZh.Range("K3").Value = "2016-11-29 02:46:02"

# Handback is in sync now, so the stale "handback not latest" error
# details are cleared.
$wsZh.Range("P2").Value = ""
$wsZh.Range("P3").Value = ""

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

# New handback datetimes recorded for this run.
$wsDe.Range("K2").Value = "2016-11-29 02:46:20"
$wsDe.Range("K3").Value = "2016-11-29 02:46:20"

$wsDe.Range("P2").Value = ""
$wsDe.Range("P3").Value = ""

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(16).ColumnWidth = 12.833333333333334
